# Fruta / hortaliza, semanal
# Insert a new weekly record at row 92 (pushing the previous rows 92-112 down
# to 93-113), re-using the descriptive fields of the (old) row 92 entry but
# with updated date / price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current contents of row 92 (columns A-T) before shifting rows
# down, so we can reuse the unchanged descriptive fields for the newly
# inserted row.
$oldValues = @{}
for ($col = 1; $col -le 20; $col++) {
    $oldValues[$col] = $ws.Cells.Item(92, $col).Value2
}

# Shift rows 92:112 down to 93:113, leaving a blank row 92 behind.
$ws.Rows.Item(92).Insert()

# Re-populate the new row 92 with the previous row's values.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(92, $col).Value = $oldValues[$col]
}

# Apply the updated figures for the newly inserted record.
$ws.Cells.Item(92, 4).Value = 44776    # D92 - Fecha
$ws.Cells.Item(92, 14).Value = 14000   # N92 - Precio mínimo
$ws.Cells.Item(92, 15).Value = 15000   # O92 - Precio máximo
$ws.Cells.Item(92, 16).Value = 14500   # P92 - Precio promedio ponderado
$ws.Cells.Item(92, 19).Value = 725     # S92 - Precio $/Kg
